# ===========================================================================
# karma_performance.xlsx edit script
#   Sheet1 "Sheet1"            -> "Single double"   (keeps D3:P10 table + chart1)
#   Sheet2 "Sheet2" (empty)    -> "Sequence of items" (gets old D36:P42 table
#                                 renumbered to D3:P9, plus chart2)
#   Sheet3 "Sheet3" (empty)    -> "Single int" (gets new int-formatting table
#                                 + brand-new chart3), becomes the active sheet
# ===========================================================================

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- 1. Rename sheets -------------------------------------------------------
$ws1.Name = "Single double"
$ws2.Name = "Sequence of items"
$ws3.Name = "Single int"

# --- 2. Relocate the "sequence of items" table from Sheet1 rows 36-42 ------
#        to Sheet2 rows 3-9 (with formatting) ------------------------------
$ws1.Range("D36:P42").Copy()
$ws2.Range("D3").PasteSpecial()
$ws1.Range("D36:P42").ClearContents()
$ws1.Rows("36:42").Delete()

# --- 3. Remove chart2 ("Chart 3") from Sheet1 - it belongs on Sheet2 now ---
$ws1ChartObjects = $ws1.ChartObjects()
$ws1ChartObjects.Item(2).Delete()

# --- 4. Fix up chart1's series references (Sheet1 -> 'Single double') ------
$chart1 = $ws1.ChartObjects().Item(1).Chart
for ($i = 1; $i -le $chart1.SeriesCollection().Count; $i++) {
    $s = $chart1.SeriesCollection().Item($i)
    $row = 3 + $i
    $s.Name = "='Single double'!`$D`$$row"
    $s.XValues = "='Single double'!`$E`$3:`$J`$3"
    $s.Values = "='Single double'!`$E`$$row`:`$J`$$row"
}

Write-Host "Steps 1-4 done"
